$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("E2").Value = "2016-04-12 06:50:20"
$ws.Range("H2").Value = "2016-04-12 16:11:00"
